# The workbook originally contained 6 data rows (rows 2-7):
#   rows 2-4: Sending cluster = ECs   (ligand C1qa / receptor Cspg4) -> target clusters ECs/FAPs/MuSCs
#   rows 5-7: Sending cluster = MuSCs (ligand C1qa / receptor Cspg4) -> target clusters ECs/FAPs/MuSCs
#
# The new TPM recomputation removes the "ECs" sending-cluster rows entirely and keeps
# only the "MuSCs" sending-cluster rows (which become rows 2-4), with refreshed
# expression/specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "ECs" sending-cluster rows (rows 2-4). This shifts the former
# "MuSCs" sending-cluster rows (5-7) up to become rows 2-4.
$ws.Rows("2:4").Delete()

function Set-RowValues($rowIndex, $values) {
    $arr = New-Object 'object[,]' 1, $values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $rangeAddr = "I" + $rowIndex + ":T" + $rowIndex
    $ws.Range($rangeAddr).Value = $arr
}

# Refresh the numeric columns (I through T) for the three remaining rows with the
# updated TPM-derived values.
Set-RowValues 2 @(1, 1, 3, 1, 1.366995666666667, 4.100987, 0.02653821474268573, 0.02653821474268573, 0.02312592135822222, 0.208133292224, 0.02653821474268573, 0.02653821474268573)
Set-RowValues 3 @(1, 1, 3, 1, 12.844987, 38.534961, 0.2493665720274216, 0.2493665720274215, 0.2173029267413333, 1.955726340672, 0.2493665720274216, 0.2493665720274215)
Set-RowValues 4 @(1, 1, 3, 1, 37.298478, 111.895434, 0.7240952132298927, 0.7240952132298926, 0.6309907851519999, 5.678917066367999, 0.7240952132298927, 0.7240952132298926)
